$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update Alemania (row 8)
$ws.Range("B8").Value = 138221
$ws.Range("C8").Value = 523
$ws.Range("E8").Value = 52323
$ws.Range("G8").Value = 46
$ws.Range("H8").Value = 4098

# Update Iran (row 11)
$ws.Range("B11").Value = 79494
$ws.Range("C11").Value = 1499
$ws.Range("D11").Value = 54064
$ws.Range("E11").Value = 20472
$ws.Range("F11").Value = 3563
$ws.Range("G11").Value = 89
$ws.Range("H11").Value = 4958

# Swap Azerbaiyan / Eslovenia rows (row 72 and 73) and update Eslovenia's stats
$ws.Range("A72").Value = "Eslovenia"
$ws.Range("B72").Value = 1304
$ws.Range("C72").Value = 36
$ws.Range("D72").Value = 174
$ws.Range("E72").Value = 1064
$ws.Range("F72").Value = 28
$ws.Range("G72").Value = 5
$ws.Range("H72").Value = 66

$ws.Range("A73").Value = "Azerbaiyan"
$ws.Range("B73").Value = 1283
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 460
$ws.Range("E73").Value = 808
$ws.Range("F73").Value = 28
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 15

# Update Malta (row 105)
$ws.Range("D105").Value = 91
$ws.Range("E105").Value = 328

# Update Senegal (row 111)
$ws.Range("B111").Value = 342
$ws.Range("C111").Value = 7
$ws.Range("D111").Value = 198
$ws.Range("E111").Value = 142
